# Adds two new columns to the netlist table ("BL" and "Operation Freq"),
# initialised to 0 for every component row, per the commit
# "Correciones con lectura de Excel para distribuidos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1, G1) -- same default style as the rest of row 1's
# un-styled cells (they carry no explicit style in the target file either).
$ws.Range("F1").Value = "BL"
$ws.Range("G1").Value = "Operation Freq"

# New data cells for each of the 4 component rows, all defaulting to 0.
$ws.Range("F2:G5").Value = 0

# Reproduce the final selection recorded in the saved workbook.
$ws.Range("C6").Select()
